$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New sample data: common names + bio-verification markers for a few rows ---
# Values are written in the exact order the original authoring tool created
# the shared-string table in, so newly-introduced strings line up with the
# canonical file (re-used strings like "Common Name" just reference the
# existing shared string).

$ws.Range("E2").Value = "N"

$ws.Range("B3").Value = "Kuranda Tree Frog"
$ws.Range("E3").Value = "Common Name"
$ws.Range("H3").Value = "Common Name?"

$ws.Range("B4").Value = "Peron's Tree Frog"
$ws.Range("H4").Value = "Common Name?"

$ws.Range("H5").Value = "Common Name"

$ws.Range("H6").Value = "Common Name"

$ws.Range("B7").Value = "Orange Thighed Tree Frog"
$ws.Range("H7").Value = "Common Name?"

$ws.Range("H8").Value = "Common Name"

$ws.Range("H9").Value = "Common Name"

$ws.Range("H10").Value = "Common Name"

$ws.Range("B11").Value = "Purple-crowned fairy wren"
$ws.Range("H11").Value = "Common Name?"

# --- Column H (Bio-verified) is now wider to fit the new sample values ---
$ws.Columns("H").ColumnWidth = 16.3

# --- Selection moved down onto the newly-populated sample rows ---
$ws.Activate()
$ws.Range("B12").Select()
